# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-391) on the active sheet from 2023-09-13 (serial 45182) to
# 2023-09-15 (serial 45184), matching the source workbook's refreshed
# "last changed" timestamp for every logging-notification record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C391").Value = 45184
